{"js": "// Replace each division-problem text in the table with its new value.\n// Mapping derived from the authoritative diff (old -> new).\nconst replacements = [\n  [\"321\u00f73=\", \"949\u00f79=\"],\n  [\"135\u00f72=\", \"877\u00f77=\"],\n  [\"182\u00f77=\", \"168\u00f79=\"],\n  [\"130\u00f77=\", \"566\u00f78=\"],\n  [\"150\u00f74=\", \"577\u00f72=\"],\n  [\"916\u00f77=\", \"824\u00f78=\"],\n  [\"622\u00f77=\", \"346\u00f76=\"],\n  [\"158\u00f75=\", \"520\u00f74=\"],\n  [\"219\u00f78=\", \"443\u00f79=\"],\n  [\"359\u00f76=\", \"791\u00f79=\"],\n  [\"434\u00f74=\", \"545\u00f74=\"],\n  [\"140\u00f79=\", \"610\u00f77=\"],\n  [\"578\u00f73=\", \"482\u00f77=\"],\n  [\"150\u00f76=\", \"712\u00f78=\"],\n  [\"819\u00f73=\", \"152\u00f74=\"],\n  [\"288\u00f72=\", \"247\u00f77=\"],\n  [\"328\u00f75=\", \"862\u00f75=\"],\n  [\"205\u00f72=\", \"363\u00f75=\"],\n  [\"239\u00f79=\", \"949\u00f72=\"],\n  [\"372\u00f73=\", \"809\u00f78=\"],\n  [\"939\u00f77=\", \"576\u00f72=\"],\n  [\"278\u00f79=\", \"172\u00f74=\"],\n  [\"946\u00f79=\", \"885\u00f77=\"],\n  [\"634\u00f75=\", \"121\u00f75=\"],\n  [\"347\u00f76=\", \"984\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each division-problem text in the table with its new value.\n# Mapping derived from the authoritative diff (old -> new).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"321\u00f73=\", \"949\u00f79=\"),\n  @(\"135\u00f72=\", \"877\u00f77=\"),\n  @(\"182\u00f77=\", \"168\u00f79=\"),\n  @(\"130\u00f77=\", \"566\u00f78=\"),\n  @(\"150\u00f74=\", \"577\u00f72=\"),\n  @(\"916\u00f77=\", \"824\u00f78=\"),\n  @(\"622\u00f77=\", \"346\u00f76=\"),\n  @(\"158\u00f75=\", \"520\u00f74=\"),\n  @(\"219\u00f78=\", \"443\u00f79=\"),\n  @(\"359\u00f76=\", \"791\u00f79=\"),\n  @(\"434\u00f74=\", \"545\u00f74=\"),\n  @(\"140\u00f79=\", \"610\u00f77=\"),\n  @(\"578\u00f73=\", \"482\u00f77=\"),\n  @(\"150\u00f76=\", \"712\u00f78=\"),\n  @(\"819\u00f73=\", \"152\u00f74=\"),\n  @(\"288\u00f72=\", \"247\u00f77=\"),\n  @(\"328\u00f75=\", \"862\u00f75=\"),\n  @(\"205\u00f72=\", \"363\u00f75=\"),\n  @(\"239\u00f79=\", \"949\u00f72=\"),\n  @(\"372\u00f73=\", \"809\u00f78=\"),\n  @(\"939\u00f77=\", \"576\u00f72=\"),\n  @(\"278\u00f79=\", \"172\u00f74=\"),\n  @(\"946\u00f79=\", \"885\u00f77=\"),\n  @(\"634\u00f75=\", \"121\u00f75=\"),\n  @(\"347\u00f76=\", \"984\u00f75=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $oldText\n  $rng.Find.Replacement.Text = $newText\n  $rng.Find.Forward = $true\n  $rng.Find.MatchCase = $true\n  $rng.Find.MatchWholeWord = $false\n  $rng.Find.MatchWildcards = $false\n  $rng.Find.Execute(\n    $oldText,\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n  )\n}\n"}
